$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldText) {
        $cell.Value = $newText
    }
}
